$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the login e-mail address with a mailto hyperlink in T2
$cellT2 = $ws.Range("T2")
$hlink = $ws.Hyperlinks.Add($cellT2, "mailto:taousautotester@hpe.com", "", "", "taousautotester@hpe.com")
$hlink.TextToDisplay = "mailto:taousautotester@hpe.com"
$cellT2.Value = "taousautotester@hpe.com"

# Update the password/token value in U2
$ws.Range("U2").Value = "5810ca086fd249fe54f03436d5829007179d176ceef6d120c899"

# Update the view: scroll/selection now centred on the new login columns
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("T1:U2").Select()
